$d = $word.ActiveDocument
$r = $d.Content
$r.Find.Execute("renunta la rezervare", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
Write-Host $r.Text
